$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Correlation")

# Insert 4 new columns before the existing "correlation" column (E),
# shifting it (and its data) to column I.
$ws.Range("E1:H1").EntireColumn.Insert()

# Add headers for the newly inserted columns.
$ws.Range("E1").Value = "Species_1"
$ws.Range("F1").Value = "Species_2"
$ws.Range("G1").Value = "Scale_1"
$ws.Range("H1").Value = "Scale_2"

# Update the active selection to match the target workbook state.
$ws.Range("H1").Select()
